$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 194.5
$ws.Range("I2").Value = 211.57143
$ws.Range("K2").Value = 211.57143
$ws.Range("M2").Value = -98.57142999999999

$ws.Range("H28").Value = 835.1429000000001
$ws.Range("I28").Value = 806.6667
$ws.Range("J28").Value = 1006
$ws.Range("K28").Value = 806.6667
$ws.Range("L28").Value = 1006
$ws.Range("M28").Value = -321.6667
$ws.Range("N28").Value = -1976

$ws.Range("H70").Value = 7262
$ws.Range("J70").Value = 9366.799999999999
$ws.Range("L70").Value = 28100.4
$ws.Range("N70").Value = -28640.4

$ws.Range("H73").Value = 7262
$ws.Range("J73").Value = 9366.799999999999
$ws.Range("L73").Value = 28100.4
$ws.Range("N73").Value = -29972.4

$ws.Range("H96").Value = 1018.1667
$ws.Range("I96").Value = 1569.6666
$ws.Range("J96").Value = 466.66666
$ws.Range("K96").Value = 4708.9998
$ws.Range("L96").Value = 1399.99998
$ws.Range("M96").Value = -3335.9998
$ws.Range("N96").Value = -4145.999980000001

$ws.Range("H135").Value = 250
$ws.Range("I135").Value = 250
$ws.Range("K135").Value = 2250
$ws.Range("M135").Value = 285

$ws.Range("H137").Value = 4572.769
$ws.Range("I137").Value = 3937.25
$ws.Range("K137").Value = 11811.75
$ws.Range("M137").Value = -9261.75

$ws.Range("H141").Value = 890
$ws.Range("I141").Value = 890
$ws.Range("K141").Value = 2670
$ws.Range("M141").Value = 2510

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2468.1875
$ws.Range("I61").Value = 2486.0667
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 2486.0667
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -2274.0667
$ws.Range("N61").Value = -2624

$ws.Range("H136").Value = 2468.1875
$ws.Range("I136").Value = 2486.0667
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 7458.2001
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -4908.2001
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 551.5
$ws.Range("I80").Value = 349.5
$ws.Range("K80").Value = 349.5
$ws.Range("M80").Value = 648.5

$ws.Range("H83").Value = 551.5
$ws.Range("I83").Value = 349.5
$ws.Range("K83").Value = 1747.5
$ws.Range("M83").Value = 3244.5

$ws.Range("H86").Value = 974.75
$ws.Range("J86").Value = 933
$ws.Range("L86").Value = 933
$ws.Range("N86").Value = -3179

$ws.Range("H88").Value = 21114
$ws.Range("J88").Value = 21114
$ws.Range("L88").Value = 21114
$ws.Range("N88").Value = -21926

$ws.Range("H89").Value = 974.75
$ws.Range("J89").Value = 933
$ws.Range("L89").Value = 4665
$ws.Range("N89").Value = -15897

$ws.Range("H91").Value = 21114
$ws.Range("J91").Value = 21114
$ws.Range("L91").Value = 21114
$ws.Range("N91").Value = -23922

$ws.Range("H94").Value = 8556
$ws.Range("I94").Value = 7867.5
$ws.Range("K94").Value = 7867.5
$ws.Range("M94").Value = -7416.5

$ws.Range("H95").Value = 3783.5
$ws.Range("J95").Value = 3783.5
$ws.Range("L95").Value = 3783.5
$ws.Range("N95").Value = -9275.5

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1247.5
$ws.Range("I32").Value = 995
$ws.Range("K32").Value = 995
$ws.Range("M32").Value = -679

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 771.25
$ws.Range("I97").Value = 771.25
$ws.Range("K97").Value = 2313.75
$ws.Range("M97").Value = -1817.75

$ws.Range("H109").Value = 227
$ws.Range("I109").Value = 227
$ws.Range("K109").Value = 681
$ws.Range("M109").Value = 359

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1449.75
$ws.Range("I113").Value = 1266.3334
$ws.Range("K113").Value = 1266.3334
$ws.Range("M113").Value = 903.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2299
$ws.Range("I32").Value = 2299
$ws.Range("K32").Value = 2299
$ws.Range("M32").Value = -1982

$ws.Range("H76").Value = 20587.3
$ws.Range("J76").Value = 20731.889
$ws.Range("L76").Value = 20731.889
$ws.Range("N76").Value = -21407.889

$ws.Range("H79").Value = 20587.3
$ws.Range("J79").Value = 20731.889
$ws.Range("L79").Value = 20731.889
$ws.Range("N79").Value = -23071.889

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I4").Value = 27219.8
$ws.Range("J4").Value = 2133.25
$ws.Range("K4").Value = 27219.8
$ws.Range("L4").Value = 2133.25
$ws.Range("M4").Value = -27106.8
$ws.Range("N4").Value = -2359.25

$ws.Range("H5").Value = 5000374.5
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H14").Value = 2004
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H19").Value = 18999
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 18999
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 18999
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -19347

$ws.Range("H40").Value = 55006.25
$ws.Range("I40").Value = 20012.5
$ws.Range("K40").Value = 20012.5
$ws.Range("M40").Value = -19863.5

$ws.Range("H68").Value = 26366.666
$ws.Range("J68").Value = 26366.666
$ws.Range("L68").Value = 26366.666
$ws.Range("N68").Value = -27988.666

$ws.Range("H71").Value = 26366.666
$ws.Range("J71").Value = 26366.666
$ws.Range("L71").Value = 79099.99800000001
$ws.Range("N71").Value = -87211.99800000001

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H136").Value = 1140.3572
$ws.Range("I136").Value = 1140.3572
$ws.Range("K136").Value = 3421.0716
$ws.Range("M136").Value = -871.0715999999998
